$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force text format on D and E columns (rows 2-51) so numeric-looking
# strings like "226.60" or "1.740.22" are preserved as text, matching
# the original inlineStr cell storage.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.959.92'
$ws.Range("E2").Value = '  -4.20%  '
$ws.Range("D3").Value = '1.740.22'
$ws.Range("E3").Value = '  -4.60%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '226.60'
$ws.Range("D6").Value = '0.5790'
$ws.Range("E6").Value = '  -3.32%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '0.2741'
$ws.Range("E8").Value = '  -0.74%  '
$ws.Range("D9").Value = '23.15'
$ws.Range("E9").Value = '  -1.33%  '
$ws.Range("D10").Value = '0.06627'
$ws.Range("E10").Value = '  -4.51%  '
$ws.Range("D11").Value = '0.07546'
$ws.Range("E11").Value = '  -0.70%  '
$ws.Range("D12").Value = '1.741.77'
$ws.Range("E12").Value = '  -4.48%  '
$ws.Range("D13").Value = '4.709'
$ws.Range("D14").Value = '0.6020'
$ws.Range("E14").Value = '  -3.93%  '
$ws.Range("D15").Value = '1.976.82'
$ws.Range("E15").Value = '  -4.60%  '
$ws.Range("D16").Value = '74.62'
$ws.Range("E16").Value = '  -3.45%  '
$ws.Range("D17").Value = '0.000008731'
$ws.Range("E17").Value = '  -10.79%  '
$ws.Range("D18").Value = '27.943.83'
$ws.Range("E18").Value = '  -3.63%  '
$ws.Range("D19").Value = '5.311'
$ws.Range("E19").Value = '  -3.82%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("D21").Value = '205.65'
$ws.Range("E21").Value = '  -4.74%  '
$ws.Range("E22").Value = '  -2.31%  '
$ws.Range("D23").Value = '6.630'
$ws.Range("E23").Value = '  -2.91%  '
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("D25").Value = '150.33'
$ws.Range("E25").Value = '  -3.35%  '
$ws.Range("D26").Value = '8.022'
$ws.Range("E26").Value = '  +1.00%  '
$ws.Range("D27").Value = '0.1234'
$ws.Range("E27").Value = '  -4.15%  '
$ws.Range("E28").Value = '  -1.91%  '
$ws.Range("D29").Value = '1.389'
$ws.Range("E29").Value = '  -2.62%  '
$ws.Range("D30").Value = '0.06186'
$ws.Range("E30").Value = '  -4.05%  '
$ws.Range("E31").Value = '  -3.25%  '
$ws.Range("E32").Value = '  -1.64%  '
$ws.Range("D33").Value = '3.737'
$ws.Range("E33").Value = '  -1.00%  '
$ws.Range("D34").Value = '1.678'
$ws.Range("E34").Value = '  -2.25%  '
$ws.Range("E35").Value = '  -4.97%  '
$ws.Range("D36").Value = '0.6394'
$ws.Range("E36").Value = '  -0.85%  '
$ws.Range("D37").Value = '2.438'
$ws.Range("E37").Value = '  -4.01%  '
$ws.Range("D38").Value = '2.717'
$ws.Range("E38").Value = '  -1.06%  '
$ws.Range("E39").Value = '  -4.32%  '
$ws.Range("D40").Value = '1.124.44'
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("D41").Value = '6.157'
$ws.Range("E41").Value = '  -6.45%  '
$ws.Range("D42").Value = '0.8769'
$ws.Range("E42").Value = '  -1.57%  '
$ws.Range("D43").Value = '1.004'
$ws.Range("E43").Value = '  +0.11%  '
$ws.Range("E44").Value = '  -0.43%  '
$ws.Range("D45").Value = '1.888.63'
$ws.Range("E45").Value = '  -4.78%  '
$ws.Range("D46").Value = '59.34'
$ws.Range("E46").Value = '  -4.33%  '
$ws.Range("D47").Value = '1.577'
$ws.Range("E47").Value = '  -2.11%  '
$ws.Range("E48").Value = '  -5.25%  '
$ws.Range("D49").Value = '8.257'
$ws.Range("E49").Value = '  -2.24%  '
$ws.Range("E50").Value = '  -2.25%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '0.4412'
$ws.Range("E51").Value = '  -2.51%  '
